# Append two new "sign up" test-data blocks (rows 32-51) to Sheet1,
# mirroring the existing blocks of 10 rows already in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row 32
$ws.Range("A32").Value = "mAlfonso"
$ws.Range("B32").Value = "MDoyle"
$ws.Range("C32").Value = "1-912-244-9624"
$ws.Range("D32").Value = "watson.leuschke@hotmail.com"
$ws.Range("E32").Value = "43043468MNon"
$ws.Range("F32").Value = "INValid first name"

# row 33
$ws.Range("A33").Value = "MStanton"
$ws.Range("B33").Value = "mWillms"
$ws.Range("C33").Value = "1-912-244-9624"
$ws.Range("D33").Value = "violet.schuster@gmail.com"
$ws.Range("E33").Value = "43043468MNon"
$ws.Range("F33").Value = "INValid last name"

# row 34
$ws.Range("A34").Value = "mAlfonso"
$ws.Range("B34").Value = "mWillms"
$ws.Range("C34").Value = "1-912-244-9624"
$ws.Range("D34").Value = "reinhold.mayer@gmail.com"
$ws.Range("E34").Value = "43043468MNon"
$ws.Range("F34").Value = "INValid first and lastname"

# row 35
$ws.Range("A35").Value = "MStanton"
$ws.Range("B35").Value = "MStanton"
$ws.Range("C35").Value = "1-912-244-9624"
$ws.Range("D35").Value = "layla.schuster@gmail.com"
$ws.Range("E35").Value = "43043468MNon"
$ws.Range("F35").Value = "INValidfirstname and last name not matched"

# row 36
$ws.Range("A36").Value = "MStanton"
$ws.Range("B36").Value = "MDoyle"
$ws.Range("C36").Value = "trst456789"
$ws.Range("D36").Value = "neoma.kassulke@yahoo.com"
$ws.Range("E36").Value = "43043468MNon"
$ws.Range("F36").Value = "INValidMobileNumber"

# row 37
$ws.Range("A37").Value = "MStanton"
$ws.Range("B37").Value = "MDoyle"
$ws.Range("C37").Value = "1-912-244-9624"
$ws.Range("D37").Value = "test@test"
$ws.Range("E37").Value = "43043468MNon"
$ws.Range("F37").Value = "INValidEmail"

# row 38
$ws.Range("A38").Value = "MStanton"
$ws.Range("B38").Value = "MDoyle"
$ws.Range("C38").Value = "1-912-244-9624"
$ws.Range("D38").Value = "lurline.jacobson@hotmail.com"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "523049"
$ws.Range("F38").Value = "INValid PasswordPassword Must be 8 characters not 6"

# row 39
$ws.Range("A39").Value = "MStanton"
$ws.Range("B39").Value = "MDoyle"
$ws.Range("C39").Value = "1-912-244-9624"
$ws.Range("D39").Value = "claudine.hammes@hotmail.com"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "523049"
$ws.Range("F39").Value = "INValid Password"

# row 40
$ws.Range("A40").Value = "MStanton"
$ws.Range("B40").Value = "MDoyle"
$ws.Range("C40").Value = "1-912-244-9624"
$ws.Range("D40").Value = "orie.bradtke@gmail.com"
$ws.Range("E40").Value = "43043468MNon"
$ws.Range("F40").Value = "Valid"

# row 41
$ws.Range("A41").Value = "MStanton"
$ws.Range("B41").Value = "MDoyle"
$ws.Range("C41").Value = "1-912-244-9624"
$ws.Range("D41").Value = "orie.bradtke@gmail.com"
$ws.Range("E41").Value = "43043468MNon"
$ws.Range("F41").Value = "INValid ThisEmailExist"

# row 42
$ws.Range("A42").Value = "mSuzanne"
$ws.Range("B42").Value = "MCronin"
$ws.Range("C42").Value = "883-943-8971"
$ws.Range("D42").Value = "orval.jakubowski@yahoo.com"
$ws.Range("E42").Value = "13635568MNon"
$ws.Range("F42").Value = "INValid first name"

# row 43
$ws.Range("A43").Value = "MNia"
$ws.Range("B43").Value = "mCorwin"
$ws.Range("C43").Value = "883-943-8971"
$ws.Range("D43").Value = "marcelo.tillman@hotmail.com"
$ws.Range("E43").Value = "13635568MNon"
$ws.Range("F43").Value = "INValid last name"

# row 44
$ws.Range("A44").Value = "mSuzanne"
$ws.Range("B44").Value = "mCorwin"
$ws.Range("C44").Value = "883-943-8971"
$ws.Range("D44").Value = "oren.kozey@yahoo.com"
$ws.Range("E44").Value = "13635568MNon"
$ws.Range("F44").Value = "INValid first and lastname"

# row 45
$ws.Range("A45").Value = "MNia"
$ws.Range("B45").Value = "MNia"
$ws.Range("C45").Value = "883-943-8971"
$ws.Range("D45").Value = "chase.adams@hotmail.com"
$ws.Range("E45").Value = "13635568MNon"
$ws.Range("F45").Value = "INValidfirstname and last name not matched"

# row 46
$ws.Range("A46").Value = "MNia"
$ws.Range("B46").Value = "MCronin"
$ws.Range("C46").Value = "trst456789"
$ws.Range("D46").Value = "hallie.ryan@hotmail.com"
$ws.Range("E46").Value = "13635568MNon"
$ws.Range("F46").Value = "INValidMobileNumber"

# row 47
$ws.Range("A47").Value = "MNia"
$ws.Range("B47").Value = "MCronin"
$ws.Range("C47").Value = "883-943-8971"
$ws.Range("D47").Value = "test@test"
$ws.Range("E47").Value = "13635568MNon"
$ws.Range("F47").Value = "INValidEmail"

# row 48
$ws.Range("A48").Value = "MNia"
$ws.Range("B48").Value = "MCronin"
$ws.Range("C48").Value = "883-943-8971"
$ws.Range("D48").Value = "gaylord.reinger@yahoo.com"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "359433"
$ws.Range("F48").Value = "INValid PasswordPassword Must be 8 characters not 6"

# row 49
$ws.Range("A49").Value = "MNia"
$ws.Range("B49").Value = "MCronin"
$ws.Range("C49").Value = "883-943-8971"
$ws.Range("D49").Value = "bianka.legros@hotmail.com"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "359433"
$ws.Range("F49").Value = "INValid Password"

# row 50
$ws.Range("A50").Value = "MNia"
$ws.Range("B50").Value = "MCronin"
$ws.Range("C50").Value = "883-943-8971"
$ws.Range("D50").Value = "letitia.cormier@gmail.com"
$ws.Range("E50").Value = "13635568MNon"
$ws.Range("F50").Value = "Valid"

# row 51
$ws.Range("A51").Value = "MNia"
$ws.Range("B51").Value = "MCronin"
$ws.Range("C51").Value = "883-943-8971"
$ws.Range("D51").Value = "letitia.cormier@gmail.com"
$ws.Range("E51").Value = "13635568MNon"
$ws.Range("F51").Value = "INValid ThisEmailExist"
